$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "41.192.96"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "2.176.05"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'254.10"
$ws.Range("E5").Value = "  +5.24%  "
$ws.Range("D6").Value = "'0.625"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("D7").Value = "'67.64"
$ws.Range("E7").Value = "  -1.15%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.573"
$ws.Range("E9").Value = "  +6.99%  "
$ws.Range("D10").Value = "'37.42"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("D11").Value = "'58.95"
$ws.Range("E11").Value = "  +2.64%  "
$ws.Range("D12").Value = "'0.0928"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").Value = "'7.06"
$ws.Range("E13").Value = "  +7.93%  "
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "2.501.17"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").Value = "'0.867"
$ws.Range("E16").Value = "  +4.77%  "
$ws.Range("D17").Value = "'14.40"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "2.184.96"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "41.126.81"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").Value = "0.0₃0949"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("D21").Value = "'6.14"
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("D22").Value = "'71.59"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").Value = "'231.38"
$ws.Range("E23").Value = "  +1.03%  "
$ws.Range("D24").Value = "'2.02"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("D25").Value = "'3.95"
$ws.Range("E25").Value = "  +10.23%  "
$ws.Range("D26").Value = "'11.74"
$ws.Range("E26").Value = "  +21.83%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +5.80%  "
$ws.Range("D30").Value = "'168.14"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").Value = "'20.58"
$ws.Range("E31").Value = "  +2.26%  "
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "'0.0747"
$ws.Range("E33").Value = "  +7.41%  "
$ws.Range("E34").Value = "  +0.23%  "
$ws.Range("D35").Value = "'5.43"
$ws.Range("E35").Value = "  +6.51%  "
$ws.Range("D36").Value = "'26.53"
$ws.Range("E36").Value = "  +13.75%  "
$ws.Range("D37").Value = "'4.61"
$ws.Range("E37").Value = "  +1.39%  "
$ws.Range("D38").Value = "'4.10"
$ws.Range("E38").Value = "  +7.38%  "
$ws.Range("E39").Value = "  +13.09%  "
$ws.Range("D40").Value = "'2.19"
$ws.Range("E40").Value = "  -2.82%  "
$ws.Range("D41").Value = "'12.50"
$ws.Range("E41").Value = "  +22.21%  "
$ws.Range("D42").Value = "'5.67"
$ws.Range("E42").Value = "  -1.69%  "
$ws.Range("D43").Value = "'64.15"
$ws.Range("E43").Value = "  +3.05%  "
$ws.Range("E44").Value = "  +5.07%  "
$ws.Range("E45").Value = "  +5.55%  "
$ws.Range("D46").Value = "'8.63"
$ws.Range("E46").Value = "  +1.04%  "
$ws.Range("D47").Value = "'0.101"
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("D48").Value = "'1.01"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("E49").Value = "  +4.50%  "
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("D51").Value = "'4.27"
$ws.Range("E51").Value = "  -4.71%  "
